# Add "2022-Q4" sheet (a copy of the "2022-Q3" detail-sheet layout/format,
# repopulated with the 2022-Q4 fund-holding data) right after the "总计"
# summary sheet, and update the "总计" summary table to include the new
# quarter (shifting the older quarters down by one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet by copying "2022-Q3" (same column
#    headers / styles as every other quarterly detail sheet) and placing
#    it immediately after "总计".
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy([System.Reflection.Missing]::Value, $wsTotal)

$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# The source sheet (2022-Q3) has 12 data rows (rows 2-13); 2022-Q4 only
# needs 4 data rows (rows 2-5), so drop the extra rows.
$wsQ4.Rows("6:13").Delete()

# Helper to write a text-typed value (mirrors the source data, which
# stores these columns as text, e.g. "014600", "0.89") without Excel's
# automatic text->number coercion, and without leaving the cell's style
# pointing at a different (quote-prefixed) format than its neighbours.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $wsQ4.Range("B2") "014600"
Set-TextValue $wsQ4.Range("C2") "博时回报严选混合A"
Set-TextValue $wsQ4.Range("D2") "0.89"
Set-TextValue $wsQ4.Range("E2") "93.75"
Set-TextValue $wsQ4.Range("F2") "7.88"
Set-TextValue $wsQ4.Range("G2") "0.0701"
$wsQ4.Range("H2").Value = 3

# Row 3
Set-TextValue $wsQ4.Range("B3") "014232"
Set-TextValue $wsQ4.Range("C3") "博时专精特新主题混合A"
Set-TextValue $wsQ4.Range("D3") "3.00"
Set-TextValue $wsQ4.Range("E3") "80.89"
Set-TextValue $wsQ4.Range("F3") "1.35"
Set-TextValue $wsQ4.Range("G3") "0.0405"
$wsQ4.Range("H3").Value = 2

# Row 4
Set-TextValue $wsQ4.Range("B4") "014233"
Set-TextValue $wsQ4.Range("C4") "博时专精特新主题混合C"
Set-TextValue $wsQ4.Range("D4") "2.58"
Set-TextValue $wsQ4.Range("E4") "80.89"
Set-TextValue $wsQ4.Range("F4") "1.35"
Set-TextValue $wsQ4.Range("G4") "0.0348"
$wsQ4.Range("H4").Value = 2

# Row 5
Set-TextValue $wsQ4.Range("B5") "014601"
Set-TextValue $wsQ4.Range("C5") "博时回报严选混合C"
Set-TextValue $wsQ4.Range("D5") "0.05"
Set-TextValue $wsQ4.Range("E5") "93.75"
Set-TextValue $wsQ4.Range("F5") "7.88"
Set-TextValue $wsQ4.Range("G5") "0.0039"
$wsQ4.Range("H5").Value = 3

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert the new 2022-Q4 row at the
#    top of the data (row 2) and push every older quarter down one row,
#    adding the newly-exposed 2021-Q2 row (row 7).
# ---------------------------------------------------------------------

# Give the new last row (A7) the same style as the row above it before
# writing its value.
$wsTotal.Range("A6").Copy()
$wsTotal.Range("A7").PasteSpecial(-4122)
$wsTotal.Application.CutCopyMode = $false

$summary = @(
    @(0, "2022-Q4", 4, 0.15),
    @(1, "2022-Q3", 12, 0.46),
    @(2, "2022-Q2", 2, 0.15),
    @(3, "2021-Q4", 15, 2.52),
    @(4, "2021-Q3", 25, 4.58),
    @(5, "2021-Q2", 1, 0.01)
)

for ($i = 0; $i -lt $summary.Length; $i++) {
    $row = $i + 2
    $values = $summary[$i]
    $wsTotal.Cells.Item($row, 1).Value = $values[0]
    $wsTotal.Cells.Item($row, 2).Value = $values[1]
    $wsTotal.Cells.Item($row, 3).Value = $values[2]
    $wsTotal.Cells.Item($row, 4).Value = $values[3]
}

# Restore the originally-active sheet (the copy/paste operations above
# leave "2022-Q4" selected).
$wsTotal.Activate()
